$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 2382.9412
$ws.Cells.Item(62, 9).Value = 2600.4443
$ws.Cells.Item(62, 10).Value = 2138.25
$ws.Cells.Item(62, 11).Value = 2600.4443
$ws.Cells.Item(62, 12).Value = 2138.25
$ws.Cells.Item(62, 13).Value = -1976.4443
$ws.Cells.Item(62, 14).Value = -3386.25

$ws.Cells.Item(65, 8).Value = 2382.9412
$ws.Cells.Item(65, 9).Value = 2600.4443
$ws.Cells.Item(65, 10).Value = 2138.25
$ws.Cells.Item(65, 11).Value = 13002.2215
$ws.Cells.Item(65, 12).Value = 10691.25
$ws.Cells.Item(65, 13).Value = -9882.2215
$ws.Cells.Item(65, 14).Value = -16931.25

$ws.Cells.Item(107, 8).Value = 649.9231
$ws.Cells.Item(107, 9).Value = 655.2917
$ws.Cells.Item(107, 10).Value = 585.5
$ws.Cells.Item(107, 11).Value = 655.2917
$ws.Cells.Item(107, 12).Value = 585.5
$ws.Cells.Item(107, 13).Value = 1264.7083
$ws.Cells.Item(107, 14).Value = -4425.5

$ws.Cells.Item(112, 8).Value = 1030.5161
$ws.Cells.Item(112, 10).Value = 1048.2
$ws.Cells.Item(112, 12).Value = 3144.6
$ws.Cells.Item(112, 14).Value = -5360.6

$ws.Cells.Item(137, 8).Value = 1169.6875
$ws.Cells.Item(137, 9).Value = 1095.08
$ws.Cells.Item(137, 10).Value = 1436.1428
$ws.Cells.Item(137, 11).Value = 3285.24
$ws.Cells.Item(137, 12).Value = 4308.428400000001
$ws.Cells.Item(137, 13).Value = -735.2399999999998
$ws.Cells.Item(137, 14).Value = -9408.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 32289.406
$ws.Cells.Item(2, 9).Value = 1099.5769
$ws.Cells.Item(2, 10).Value = 167445.33
$ws.Cells.Item(2, 11).Value = 1099.5769
$ws.Cells.Item(2, 12).Value = 167445.33
$ws.Cells.Item(2, 13).Value = -986.5769
$ws.Cells.Item(2, 14).Value = -167671.33

$ws.Cells.Item(32, 8).Value = 4605.85
$ws.Cells.Item(32, 9).Value = 4284.3813
$ws.Cells.Item(32, 10).Value = 15000
$ws.Cells.Item(32, 11).Value = 4284.3813
$ws.Cells.Item(32, 12).Value = 15000
$ws.Cells.Item(32, 13).Value = -3997.3813
$ws.Cells.Item(32, 14).Value = -15574

$ws.Cells.Item(45, 8).Value = 49594.363
$ws.Cells.Item(45, 9).Value = 68504.13
$ws.Cells.Item(45, 10).Value = 9073.429
$ws.Cells.Item(45, 11).Value = 68504.13
$ws.Cells.Item(45, 12).Value = 9073.429
$ws.Cells.Item(45, 13).Value = -68127.13
$ws.Cells.Item(45, 14).Value = -9827.429

$ws.Cells.Item(55, 8).Value = 14133.333
$ws.Cells.Item(55, 9).Value = 8000
$ws.Cells.Item(55, 10).Value = 15885.714
$ws.Cells.Item(55, 11).Value = 8000
$ws.Cells.Item(55, 12).Value = 15885.714
$ws.Cells.Item(55, 13).Value = -7685
$ws.Cells.Item(55, 14).Value = -16515.714

$ws.Cells.Item(61, 8).Value = 1762.7675
$ws.Cells.Item(61, 9).Value = 865.4167
$ws.Cells.Item(61, 10).Value = 2896.2632
$ws.Cells.Item(61, 11).Value = 865.4167
$ws.Cells.Item(61, 12).Value = 2896.2632
$ws.Cells.Item(61, 13).Value = -653.4167
$ws.Cells.Item(61, 14).Value = -3320.2632

$ws.Cells.Item(74, 8).Value = 920.1212
$ws.Cells.Item(74, 9).Value = 859.44446
$ws.Cells.Item(74, 10).Value = 1193.1666
$ws.Cells.Item(74, 11).Value = 859.44446
$ws.Cells.Item(74, 12).Value = 1193.1666
$ws.Cells.Item(74, 13).Value = 14.55553999999995
$ws.Cells.Item(74, 14).Value = -2941.1666

$ws.Cells.Item(77, 8).Value = 920.1212
$ws.Cells.Item(77, 9).Value = 859.44446
$ws.Cells.Item(77, 10).Value = 1193.1666
$ws.Cells.Item(77, 11).Value = 4297.2223
$ws.Cells.Item(77, 12).Value = 5965.833000000001
$ws.Cells.Item(77, 13).Value = 70.77769999999964
$ws.Cells.Item(77, 14).Value = -14701.833

$ws.Cells.Item(111, 8).Value = 34000
$ws.Cells.Item(111, 10).Value = 34000
$ws.Cells.Item(111, 12).Value = 34000
$ws.Cells.Item(111, 14).Value = -42180

$ws.Cells.Item(112, 8).Value = 16000
$ws.Cells.Item(112, 10).Value = 16000
$ws.Cells.Item(112, 12).Value = 16000
$ws.Cells.Item(112, 14).Value = -18954

$ws.Cells.Item(116, 8).Value = 32289.406
$ws.Cells.Item(116, 9).Value = 1099.5769
$ws.Cells.Item(116, 10).Value = 167445.33
$ws.Cells.Item(116, 11).Value = 1099.5769
$ws.Cells.Item(116, 12).Value = 167445.33
$ws.Cells.Item(116, 13).Value = 1194.4231
$ws.Cells.Item(116, 14).Value = -172033.33

$ws.Cells.Item(122, 8).Value = 2079.55
$ws.Cells.Item(122, 9).Value = 1798.1333
$ws.Cells.Item(122, 10).Value = 2923.8
$ws.Cells.Item(122, 11).Value = 5394.3999
$ws.Cells.Item(122, 12).Value = 8771.400000000001
$ws.Cells.Item(122, 13).Value = -2944.3999
$ws.Cells.Item(122, 14).Value = -13671.4

$ws.Cells.Item(132, 8).Value = 20584.156
$ws.Cells.Item(132, 9).Value = 25866.25
$ws.Cells.Item(132, 10).Value = 4737.875
$ws.Cells.Item(132, 11).Value = 77598.75
$ws.Cells.Item(132, 12).Value = 14213.625
$ws.Cells.Item(132, 13).Value = -75068.75
$ws.Cells.Item(132, 14).Value = -19273.625

$ws.Cells.Item(136, 8).Value = 1762.7675
$ws.Cells.Item(136, 9).Value = 865.4167
$ws.Cells.Item(136, 10).Value = 2896.2632
$ws.Cells.Item(136, 11).Value = 2596.2501
$ws.Cells.Item(136, 12).Value = 8688.7896
$ws.Cells.Item(136, 13).Value = -46.2501000000002
$ws.Cells.Item(136, 14).Value = -13788.7896

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 32289.406
$ws.Cells.Item(3, 9).Value = 1099.5769
$ws.Cells.Item(3, 10).Value = 167445.33
$ws.Cells.Item(3, 11).Value = 1099.5769
$ws.Cells.Item(3, 12).Value = 167445.33
$ws.Cells.Item(3, 13).Value = -985.5769
$ws.Cells.Item(3, 14).Value = -167673.33

$ws.Cells.Item(64, 8).Value = 1567.4286
$ws.Cells.Item(64, 10).Value = 1819.625
$ws.Cells.Item(64, 12).Value = 1819.625
$ws.Cells.Item(64, 14).Value = -2269.625

$ws.Cells.Item(67, 8).Value = 1567.4286
$ws.Cells.Item(67, 10).Value = 1819.625
$ws.Cells.Item(67, 12).Value = 1819.625
$ws.Cells.Item(67, 14).Value = -3379.625

$ws.Cells.Item(134, 8).Value = 13075.296
$ws.Cells.Item(134, 9).Value = 14443
$ws.Cells.Item(134, 10).Value = 4413.1665
$ws.Cells.Item(134, 11).Value = 43329
$ws.Cells.Item(134, 12).Value = 13239.4995
$ws.Cells.Item(134, 13).Value = -40794
$ws.Cells.Item(134, 14).Value = -18309.4995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 2455.5557
$ws.Cells.Item(62, 9).Value = 2150
$ws.Cells.Item(62, 11).Value = 2150
$ws.Cells.Item(62, 13).Value = -1526

$ws.Cells.Item(65, 8).Value = 2455.5557
$ws.Cells.Item(65, 9).Value = 2150
$ws.Cells.Item(65, 11).Value = 10750
$ws.Cells.Item(65, 13).Value = -7630

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 178249.14
$ws.Cells.Item(129, 9).Value = 8548
$ws.Cells.Item(129, 10).Value = 232244.95
$ws.Cells.Item(129, 11).Value = 25644
$ws.Cells.Item(129, 12).Value = 696734.8500000001
$ws.Cells.Item(129, 13).Value = -20644
$ws.Cells.Item(129, 14).Value = -706734.8500000001

$ws.Cells.Item(131, 8).Value = 892.25
$ws.Cells.Item(131, 10).Value = 896.2121
$ws.Cells.Item(131, 12).Value = 2688.6363
$ws.Cells.Item(131, 14).Value = -12768.6363

$ws.Cells.Item(134, 8).Value = 3971.7742
$ws.Cells.Item(134, 9).Value = 2732
$ws.Cells.Item(134, 10).Value = 4562.143
$ws.Cells.Item(134, 11).Value = 8196
$ws.Cells.Item(134, 12).Value = 13686.429
$ws.Cells.Item(134, 13).Value = -3126
$ws.Cells.Item(134, 14).Value = -23826.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 404301.72
$ws.Cells.Item(102, 9).Value = 3456.6667
$ws.Cells.Item(102, 10).Value = 671531.75
$ws.Cells.Item(102, 11).Value = 3456.6667
$ws.Cells.Item(102, 12).Value = 671531.75
$ws.Cells.Item(102, 13).Value = -1834.6667
$ws.Cells.Item(102, 14).Value = -674775.75

$ws.Cells.Item(111, 8).Value = 47382.168
$ws.Cells.Item(111, 10).Value = 47382.168
$ws.Cells.Item(111, 12).Value = 47382.168
$ws.Cells.Item(111, 14).Value = -53516.168

$ws.Cells.Item(113, 8).Value = 2448.4546
$ws.Cells.Item(113, 9).Value = 2641.625
$ws.Cells.Item(113, 10).Value = 1933.3334
$ws.Cells.Item(113, 11).Value = 2641.625
$ws.Cells.Item(113, 12).Value = 1933.3334
$ws.Cells.Item(113, 13).Value = -471.625
$ws.Cells.Item(113, 14).Value = -6273.3334

$ws.Cells.Item(121, 8).Value = 23149.5
$ws.Cells.Item(121, 10).Value = 23149.5
$ws.Cells.Item(121, 12).Value = 23149.5
$ws.Cells.Item(121, 14).Value = -26643.5

$ws.Cells.Item(122, 8).Value = 1167.1305
$ws.Cells.Item(122, 9).Value = 1412.5
$ws.Cells.Item(122, 10).Value = 606.2857
$ws.Cells.Item(122, 11).Value = 4237.5
$ws.Cells.Item(122, 12).Value = 1818.8571
$ws.Cells.Item(122, 13).Value = -1787.5
$ws.Cells.Item(122, 14).Value = -6718.8571

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2286.1428
$ws.Cells.Item(7, 9).Value = 1345.5454
$ws.Cells.Item(7, 10).Value = 5735
$ws.Cells.Item(7, 11).Value = 1345.5454
$ws.Cells.Item(7, 12).Value = 5735
$ws.Cells.Item(7, 13).Value = -1233.5454
$ws.Cells.Item(7, 14).Value = -5959

$ws.Cells.Item(40, 8).Value = 45774.74
$ws.Cells.Item(40, 9).Value = 112802.78
$ws.Cells.Item(40, 10).Value = 2685.2856
$ws.Cells.Item(40, 11).Value = 112802.78
$ws.Cells.Item(40, 12).Value = 2685.2856
$ws.Cells.Item(40, 13).Value = -112666.78
$ws.Cells.Item(40, 14).Value = -2957.2856

$ws.Cells.Item(126, 8).Value = 2286.1428
$ws.Cells.Item(126, 9).Value = 1345.5454
$ws.Cells.Item(126, 10).Value = 5735
$ws.Cells.Item(126, 11).Value = 4036.6362
$ws.Cells.Item(126, 12).Value = 17205
$ws.Cells.Item(126, 13).Value = -1566.6362
$ws.Cells.Item(126, 14).Value = -22145

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(21, 8).Value = 21259.572
$ws.Cells.Item(21, 10).Value = 21259.572
$ws.Cells.Item(21, 12).Value = 21259.572
$ws.Cells.Item(21, 14).Value = -21729.572

$ws.Cells.Item(35, 8).Value = 21259.572
$ws.Cells.Item(35, 10).Value = 21259.572
$ws.Cells.Item(35, 12).Value = 21259.572
$ws.Cells.Item(35, 14).Value = -21839.572

$ws.Cells.Item(62, 8).Value = 5496796.5
$ws.Cells.Item(62, 9).Value = 10990893
$ws.Cells.Item(62, 11).Value = 10990893
$ws.Cells.Item(62, 13).Value = -10990269

$ws.Cells.Item(65, 8).Value = 5496796.5
$ws.Cells.Item(65, 9).Value = 10990893
$ws.Cells.Item(65, 11).Value = 54954465
$ws.Cells.Item(65, 13).Value = -54951345

$ws.Cells.Item(126, 8).Value = 1252.36
$ws.Cells.Item(126, 9).Value = 1286.7894
$ws.Cells.Item(126, 10).Value = 1143.3334
$ws.Cells.Item(126, 11).Value = 3860.3682
$ws.Cells.Item(126, 12).Value = 3430.0002
$ws.Cells.Item(126, 13).Value = -1390.3682
$ws.Cells.Item(126, 14).Value = -8370.0002

$ws.Cells.Item(136, 8).Value = 1556.8644
$ws.Cells.Item(136, 9).Value = 553.36664
$ws.Cells.Item(136, 10).Value = 2594.9656
$ws.Cells.Item(136, 11).Value = 1660.09992
$ws.Cells.Item(136, 12).Value = 7784.8968
$ws.Cells.Item(136, 13).Value = 889.9000800000001
$ws.Cells.Item(136, 14).Value = -12884.8968
